# Daily attendance processing - 2026-02-14 11:53:53 UTC
# Swap the "Recorded By" name/role ordering in column G:
#   "Miss Dina Nasr, Administrator" -> "Administrator, Miss Dina Nasr"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "Miss Dina Nasr, Administrator"
$newValue = "Administrator, Miss Dina Nasr"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($i = $firstRow; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
